$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 202.82608
$ws.Range("I33").Value = 146.11111
$ws.Range("K33").Value = 146.11111
$ws.Range("M33").Value = 82.88889
$ws.Range("H43").Value = 3315.3333
$ws.Range("I43").Value = 3598
$ws.Range("J43").Value = 2750
$ws.Range("K43").Value = 3598
$ws.Range("L43").Value = 2750
$ws.Range("M43").Value = -3529
$ws.Range("N43").Value = -2888
$ws.Range("H51").Value = 6750
$ws.Range("J51").Value = 6750
$ws.Range("L51").Value = 6750
$ws.Range("N51").Value = -7718
$ws.Range("H55").Value = 319.58823
$ws.Range("J55").Value = 410.22223
$ws.Range("L55").Value = 410.22223
$ws.Range("N55").Value = -838.2222300000001
$ws.Range("H76").Value = 4630.067
$ws.Range("I76").Value = 2950.2222
$ws.Range("J76").Value = 7149.8335
$ws.Range("K76").Value = 2950.2222
$ws.Range("L76").Value = 7149.8335
$ws.Range("M76").Value = -2635.2222
$ws.Range("N76").Value = -7779.8335
$ws.Range("H79").Value = 4630.067
$ws.Range("I79").Value = 2950.2222
$ws.Range("J79").Value = 7149.8335
$ws.Range("K79").Value = 2950.2222
$ws.Range("L79").Value = 7149.8335
$ws.Range("M79").Value = -1858.2222
$ws.Range("N79").Value = -9333.833500000001
$ws.Range("H125").Value = 2696.1
$ws.Range("J125").Value = 1827
$ws.Range("L125").Value = 16443
$ws.Range("N125").Value = -21363
$ws.Range("H127").Value = 2170.111
$ws.Range("I127").Value = 2254
$ws.Range("J127").Value = 1499
$ws.Range("K127").Value = 6762
$ws.Range("L127").Value = 4497
$ws.Range("M127").Value = -1802
$ws.Range("N127").Value = -14417

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3901.3125
$ws.Range("I45").Value = 2492.6667
$ws.Range("J45").Value = 5712.4287
$ws.Range("K45").Value = 2492.6667
$ws.Range("L45").Value = 5712.4287
$ws.Range("M45").Value = -2115.6667
$ws.Range("N45").Value = -6466.4287
$ws.Range("H61").Value = 4763674.5
$ws.Range("I61").Value = 5556787
$ws.Range("K61").Value = 5556787
$ws.Range("M61").Value = -5556575
$ws.Range("H102").Value = 16310.792
$ws.Range("I102").Value = 18308.334
$ws.Range("K102").Value = 18308.334
$ws.Range("M102").Value = -16686.334
$ws.Range("H122").Value = 2778
$ws.Range("I122").Value = 2235.0557
$ws.Range("J122").Value = 3999.625
$ws.Range("K122").Value = 6705.1671
$ws.Range("L122").Value = 11998.875
$ws.Range("M122").Value = -4255.1671
$ws.Range("N122").Value = -16898.875
$ws.Range("H132").Value = 1669441.2
$ws.Range("I132").Value = 1669441.2
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5008323.6
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5005793.6
$ws.Range("H136").Value = 4763674.5
$ws.Range("I136").Value = 5556787
$ws.Range("K136").Value = 16670361
$ws.Range("M136").Value = -16667811

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 89734.71000000001
$ws.Range("I99").Value = 104626.4
$ws.Range("J99").Value = 52505.5
$ws.Range("K99").Value = 104626.4
$ws.Range("L99").Value = 52505.5
$ws.Range("M99").Value = -103128.4
$ws.Range("N99").Value = -55501.5
$ws.Range("H134").Value = 1354202.8
$ws.Range("I134").Value = 1589613.5
$ws.Range("J134").Value = 765676
$ws.Range("K134").Value = 4768840.5
$ws.Range("L134").Value = 2297028
$ws.Range("M134").Value = -4766305.5
$ws.Range("N134").Value = -2302098

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 62998
$ws.Range("J20").Value = 62998
$ws.Range("L20").Value = 62998
$ws.Range("N20").Value = -63470
$ws.Range("H22").Value = 359.8
$ws.Range("I22").Value = 314.57144
$ws.Range("J22").Value = 465.33334
$ws.Range("K22").Value = 314.57144
$ws.Range("L22").Value = 465.33334
$ws.Range("M22").Value = 35.42856
$ws.Range("N22").Value = -1165.33334
$ws.Range("H30").Value = 62998
$ws.Range("J30").Value = 62998
$ws.Range("L30").Value = 62998
$ws.Range("N30").Value = -63180
$ws.Range("H31").Value = 48040.6
$ws.Range("I31").Value = 30695.334
$ws.Range("J31").Value = 55474.285
$ws.Range("K31").Value = 30695.334
$ws.Range("L31").Value = 55474.285
$ws.Range("M31").Value = -30400.334
$ws.Range("N31").Value = -56064.285
$ws.Range("H34").Value = 48040.6
$ws.Range("I34").Value = 30695.334
$ws.Range("J34").Value = 55474.285
$ws.Range("K34").Value = 30695.334
$ws.Range("L34").Value = 55474.285
$ws.Range("M34").Value = -30493.334
$ws.Range("N34").Value = -55878.285
$ws.Range("H86").Value = 340249.16
$ws.Range("I86").Value = 7181.3335
$ws.Range("K86").Value = 7181.3335
$ws.Range("M86").Value = -6058.3335
$ws.Range("H89").Value = 340249.16
$ws.Range("I89").Value = 7181.3335
$ws.Range("K89").Value = 35906.6675
$ws.Range("M89").Value = -30290.6675
$ws.Range("H107").Value = 837.37933
$ws.Range("I107").Value = 718.4666999999999
$ws.Range("J107").Value = 964.7857
$ws.Range("K107").Value = 718.4666999999999
$ws.Range("L107").Value = 964.7857
$ws.Range("M107").Value = 1201.5333
$ws.Range("N107").Value = -4804.7857
$ws.Range("H128").Value = 62998
$ws.Range("J128").Value = 62998
$ws.Range("L128").Value = 62998
$ws.Range("N128").Value = -72958
$ws.Range("H134").Value = 15321.1
$ws.Range("I134").Value = 20423
$ws.Range("K134").Value = 61269
$ws.Range("M134").Value = -58734

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 707.5
$ws.Range("I21").Value = 10
$ws.Range("K21").Value = 30
$ws.Range("M21").Value = 143
$ws.Range("H26").Value = 637.5
$ws.Range("I26").Value = 633
$ws.Range("K26").Value = 1899
$ws.Range("M26").Value = -1611
$ws.Range("H86").Value = 748.8333
$ws.Range("J86").Value = 624.5
$ws.Range("L86").Value = 1873.5
$ws.Range("N86").Value = -4245.5
$ws.Range("H89").Value = 748.8333
$ws.Range("J89").Value = 624.5
$ws.Range("L89").Value = 5620.5
$ws.Range("N89").Value = -17476.5
$ws.Range("H98").Value = 1373.75
$ws.Range("J98").Value = 1498.5
$ws.Range("L98").Value = 4495.5
$ws.Range("N98").Value = -7491.5
$ws.Range("H129").Value = 2178.9443
$ws.Range("J129").Value = 2872.3333
$ws.Range("L129").Value = 8616.999899999999
$ws.Range("N129").Value = -18616.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 345701.8
$ws.Range("I80").Value = 430876.25
$ws.Range("K80").Value = 430876.25
$ws.Range("M80").Value = -429878.25
$ws.Range("H83").Value = 345701.8
$ws.Range("I83").Value = 430876.25
$ws.Range("K83").Value = 2154381.25
$ws.Range("M83").Value = -2149389.25
$ws.Range("H122").Value = 6361.857
$ws.Range("J122").Value = 7124.875
$ws.Range("L122").Value = 21374.625
$ws.Range("N122").Value = -26274.625
$ws.Range("H132").Value = 67474970
$ws.Range("I132").Value = 84341224
$ws.Range("K132").Value = 253023672
$ws.Range("M132").Value = -253021142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2934.4614
$ws.Range("I40").Value = 2934.4614
$ws.Range("K40").Value = 2934.4614
$ws.Range("M40").Value = -2798.4614
$ws.Range("H68").Value = 2800
$ws.Range("I68").Value = 1600
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 1600
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -851
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 2800
$ws.Range("I71").Value = 1600
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 8000
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -4256
$ws.Range("N71").Value = -27488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11467940
$ws.Range("I132").Value = 13839143
$ws.Range("K132").Value = 41517429
$ws.Range("M132").Value = -41514899

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N132").ClearContents()
